$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new "Syarat & ketentuan" / "Isinya apa?" requirement row (row 41,
# columns B and C) - mirrors the existing rows 34-40 layout/styles.
$ws.Range("B41").Value = "Syarat & ketentuan"
$ws.Range("C41").Value = "Isinya apa?"

# The new row's text wraps onto two lines, so the row grows taller (30pt),
# same as row 6 which holds similarly-sized wrapped text.
$ws.Rows.Item(41).RowHeight = 30

# Update the view: the window had scrolled down and the selection moved to
# the newly-edited cell (C42).
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C42").Select()
